$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$ws.Range("D2").Value = '27.268.83'
$ws.Range("E2").Value = '  +0.47%  '
$ws.Range("D3").Value = '1.775.09'
$ws.Range("E3").Value = '  +3.66%  '
$ws.Range("E4").Value = '  +0.00%  '
Set-TextValue $ws.Range("D5") '313.06'
$ws.Range("E5").Value = '  +1.35%  '
Set-TextValue $ws.Range("D6") '1.000'
$ws.Range("E6").Value = '  -0.04%  '
Set-TextValue $ws.Range("D7") '0.5210'
$ws.Range("E7").Value = '  +8.86%  '
Set-TextValue $ws.Range("D8") '0.3676'
$ws.Range("E8").Value = '  +6.58%  '
Set-TextValue $ws.Range("D9") '42.87'
$ws.Range("E9").Value = '  +1.73%  '
Set-TextValue $ws.Range("D10") '0.07382'
$ws.Range("E10").Value = '  +1.36%  '
Set-TextValue $ws.Range("D11") '1.091'
$ws.Range("E11").Value = '  +4.48%  '
Set-TextValue $ws.Range("D12") '1.000'
$ws.Range("E12").Value = '  -0.04%  '
Set-TextValue $ws.Range("D13") '20.51'
$ws.Range("E13").Value = '  +3.29%  '
Set-TextValue $ws.Range("D14") '6.074'
$ws.Range("E14").Value = '  +3.59%  '
$ws.Range("D15").Value = '1.766.73'
$ws.Range("E15").Value = '  +3.28%  '
Set-TextValue $ws.Range("D16") '6.951'
$ws.Range("E16").Value = '  +1.51%  '
Set-TextValue $ws.Range("D17") '88.93'
$ws.Range("E17").Value = '  +0.04%  '
Set-TextValue $ws.Range("D18") '0.00001047'
$ws.Range("E18").Value = '  +0.38%  '
Set-TextValue $ws.Range("D19") '0.06439'
$ws.Range("E19").Value = '  +1.13%  '
Set-TextValue $ws.Range("D20") '1.000'
$ws.Range("E20").Value = '  -0.01%  '
Set-TextValue $ws.Range("D21") '16.75'
Set-TextValue $ws.Range("D22") '5.820'
$ws.Range("E22").Value = '  +3.60%  '
$ws.Range("D23").Value = '27.310.48'
$ws.Range("E23").Value = '  +0.51%  '
Set-TextValue $ws.Range("D24") '11.25'
$ws.Range("E24").Value = '  +3.89%  '
Set-TextValue $ws.Range("D25") '2.118'
$ws.Range("E25").Value = '  +1.41%  '
Set-TextValue $ws.Range("D26") '155.09'
$ws.Range("E26").Value = '  +2.11%  '
Set-TextValue $ws.Range("D27") '20.16'
$ws.Range("E27").Value = '  +2.45%  '
$ws.Range("D28").Value = '1.973.23'
$ws.Range("E28").Value = '  +3.51%  '
Set-TextValue $ws.Range("D29") '2.326'
$ws.Range("E29").Value = '  +11.27%  '
Set-TextValue $ws.Range("D30") '121.31'
$ws.Range("E30").Value = '  +1.01%  '
Set-TextValue $ws.Range("D31") '1.061'
$ws.Range("E31").Value = '  +4.39%  '
Set-TextValue $ws.Range("D32") '0.09778'
$ws.Range("E32").Value = '  +5.26%  '
Set-TextValue $ws.Range("D33") '5.567'
$ws.Range("E33").Value = '  +4.81%  '
Set-TextValue $ws.Range("D34") '3.619'
$ws.Range("E34").Value = '  +0.93%  '
Set-TextValue $ws.Range("D35") '0.02236'
$ws.Range("E35").Value = '  +1.59%  '
Set-TextValue $ws.Range("D36") '0.05967'
$ws.Range("E36").Value = '  +1.15%  '
Set-TextValue $ws.Range("D37") '11.24'
$ws.Range("E37").Value = '  +1.55%  '
Set-TextValue $ws.Range("D38") '4.836'
$ws.Range("E38").Value = '  +1.80%  '
Set-TextValue $ws.Range("D39") '0.6135'
$ws.Range("E39").Value = '  +3.32%  '
Set-TextValue $ws.Range("D40") '0.2018'
$ws.Range("E40").Value = '  +0.56%  '
Set-TextValue $ws.Range("D41") '1.433'
$ws.Range("E41").Value = '  +1.24%  '
Set-TextValue $ws.Range("D42") '8.086'
$ws.Range("E42").Value = '  +8.14%  '
Set-TextValue $ws.Range("D43") '1.139'
$ws.Range("E43").Value = '  +2.66%  '
Set-TextValue $ws.Range("D44") '13.07'
$ws.Range("E44").Value = '  +2.59%  '
Set-TextValue $ws.Range("D45") '0.5767'
$ws.Range("E45").Value = '  +2.58%  '
Set-TextValue $ws.Range("D46") '3.623'
$ws.Range("E46").Value = '  +1.44%  '
Set-TextValue $ws.Range("D47") '121.20'
$ws.Range("E47").Value = '  +2.24%  '
Set-TextValue $ws.Range("D48") '1.883'
$ws.Range("E48").Value = '  +2.27%  '
Set-TextValue $ws.Range("D49") '1.115'
$ws.Range("E49").Value = '  +2.74%  '
Set-TextValue $ws.Range("D50") '0.06708'
$ws.Range("E50").Value = '  +1.05%  '

# Row 51: coin changed from PaxDollar to Aave
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range("D51") '70.49'
$ws.Range("E51").Value = '  +1.28%  '
